# Actualizacion automatica 2025-07-21 14:40:09
# Insert a new advisor row ("FERNANDEZ MEZA JONATHAN ALEXIS") right before the
# "LOZANO MOLINA TITO JERSON" row (row 7) in both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets. The new row carries all-zero sales figures, and the
# existing rows below shift down by one. The "x de N" completion counters on
# the "VENTAS POR GRUPO" totals row move down with the insert and their
# denominator is bumped from 11 to 12 (the numerators are unchanged, since
# the inserted row contributes zero to every column).

$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO --------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(7).Insert()
$ws1.Range("A7").Value = "OFICINA-CATAECSA"
$ws1.Range("B7").Value = "FERNANDEZ MEZA JONATHAN ALEXIS"
$ws1.Range("C7:R7").Value = 0

# Totals row (was row 13, now row 14): bump "x de 11" -> "x de 12".
$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $cell = $ws1.Range($col + "14")
    $old = $cell.Value()
    $new = $old -replace "de 11", "de 12"
    $cell.Value = $new
}

# ---- Sheet: VENTA MENSUAL -------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(7).Insert()
$ws2.Range("A7").Value = "OFICINA-CATAECSA"
$ws2.Range("B7").Value = "FERNANDEZ MEZA JONATHAN ALEXIS"
$ws2.Range("C7:G7").Value = 0
# Row 14 totals on this sheet are plain numeric sums, unaffected by an
# all-zero inserted row, so nothing further to do there.
